# Generate Report for Handoff
#
# A fresh handback-report generation replaced the old source-file
# identifier (92a842dc-69ab-4d9e-8089-a90d33adc857) with a new one
# (f0b04ec4-395c-40e4-960c-294568c5e9a9) everywhere it is used -
# filenames, paths, generated xliff names - and bumped the handoff /
# generation timestamps that the run produced.

$wb = $excel.ActiveWorkbook

$oldGuid = "92a842dc-69ab-4d9e-8089-a90d33adc857"
$newGuid = "f0b04ec4-395c-40e4-960c-294568c5e9a9"

$oldHash = "14c92626108b9f57a43064cdbacc6b93a15719a0"
$newHash = "f0b6188875aeb878a50389abb76e2b525c249892"

$oldGenerateDate = "2016-09-07 06:08:37"
$newGenerateDate = "2016-09-07 06:09:24"

$oldZhHandoffDate = "2016-09-07 06:08:26"
$newZhHandoffDate = "2016-09-07 06:09:13"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/770d5626287ba534035eaa245e6f13a5087bf940/e2e/$oldGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Bulk text substitutions (keeps existing cell styles / shared-string
# reuse intact - plain Find & Replace over each sheet's used range).
# ---------------------------------------------------------------------
foreach ($ws in @($wsOverview, $wsZhCn, $wsDeDe)) {
    $ws.Cells.Replace($oldGuid, $newGuid)
    $ws.Cells.Replace($oldHash, $newHash)
}

$wsOverview.Cells.Replace($oldGenerateDate, $newGenerateDate)
$wsZhCn.Cells.Replace($oldZhHandoffDate, $newZhHandoffDate)
$wsDeDe.Cells.Replace($oldGenerateDate, $newGenerateDate)

# ---------------------------------------------------------------------
# Refresh each sheet's hyperlink so its displayed text matches the new
# path/filename while keeping the same target address.
# ---------------------------------------------------------------------
$newPathAndName = "e2e\$newGuid.md"
$newFileName = "$newGuid.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", $newPathAndName)

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, "", "", $newFileName)

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, "", "", $newFileName)

Write-Output "Report regenerated for handoff."
